$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("K3").Value = "Handback file name: kb0vmwpm.mow is different with handoff file name: 841e653b-26ff-4168-89ca-0ca720d89b40.10d715875a4c1e634675661c5c77a06f20742d8e.zh-cn."

$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("K3").Value = "Handback file name: kb0vmwpm.mow is different with handoff file name: 841e653b-26ff-4168-89ca-0ca720d89b40.10d715875a4c1e634675661c5c77a06f20742d8e.de-de."
